$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update task rows with new data (Attribution column changed from "null" to "Sanchez Rick",
# and task/status labels were reshuffled per the fix for the null-attribution bug)
$ws.Range("A2").Value = "Discover"
$ws.Range("B2").Value = "Sanchez Rick"
$ws.Range("C2").Value = "Bloqué"

$ws.Range("A3").Value = "Discover"
$ws.Range("B3").Value = "Sanchez Rick"
$ws.Range("C3").Value = "Terminé"

$ws.Range("A4").Value = "Label 1"
$ws.Range("B4").Value = "Sanchez Rick"
$ws.Range("C4").Value = "En cours"

$ws.Range("A5").Value = "Se préinscrire"
$ws.Range("B5").Value = "Sanchez Rick"
$ws.Range("C5").Value = "En cours"

$ws.Range("A6").Value = "Se préinscrire"
$ws.Range("B6").Value = "Sanchez Rick"
$ws.Range("C6").Value = "En cours"

# Column B needs to grow to fit the new, wider "Sanchez Rick" values
$ws.Columns.Item(2).ColumnWidth = 11
